$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: TRADING_ATTEMPT for ENA
$ws.Range("A22").Value = "2025-09-20T01:23:07.266452"
$ws.Range("B22").Value = "TRADING_ATTEMPT"
$ws.Range("C22").Value = "ENA"
$ws.Range("D22").Value = "UNKNOWN"
$ws.Range("E22").Value = 0.6727811902747289
$ws.Range("F22").Value = ""
$ws.Range("G22").Value = ""
$ws.Range("H22").Value = ""
$ws.Range("I22").Value = ""
$ws.Range("J22").Value = ""
$ws.Range("K22").Value = "ATTEMPT"
$ws.Range("L22").Value = "Attempting trade 1/1"

# Row 23: POSITION_OPENED for ENA
$ws.Range("A23").Value = "2025-09-20T01:23:08.787988"
$ws.Range("B23").Value = "POSITION_OPENED"
$ws.Range("C23").Value = "ENA"
$ws.Range("D23").Value = "UNKNOWN"
$ws.Range("E23").Value = 0.6727811902747289
$ws.Range("F23").Value = 2400
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 0.6326121610267794
$ws.Range("I23").Value = ""
$ws.Range("J23").Value = ""
$ws.Range("K23").Value = "SUCCESS"
$ws.Range("L23").Value = ""
